$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Update row 8 (22uH Inductor) to the new part number / datasheet link.
#    Part description ("22uH Inductor") and quantity (2) stay the same.
# ---------------------------------------------------------------------------
$ws.Range("A8").Value = "2300HT-220-H-RC"

# ---------------------------------------------------------------------------
# 2. Clear every existing hyperlink on the sheet (the runtime's
#    Hyperlinks.Delete() operates sheet-wide) and re-create all of them -
#    the original nine plus the ten new BOM rows - so relationship ids come
#    out clean with no stale/duplicate entries.
# ---------------------------------------------------------------------------
$ws.Range("F2").Hyperlinks.Delete() | Out-Null

function Add-BomLink($cellRef, $url, $display) {
    if ($display) {
        $ws.Hyperlinks.Add($ws.Range($cellRef), $url, "", "", $display) | Out-Null
    } else {
        $ws.Hyperlinks.Add($ws.Range($cellRef), $url) | Out-Null
    }
    $ws.Range($cellRef).Style = "Hyperlink"
}

$f3display = "https://www.digikey.com/en/products/detail/texas-instruments/LM2678S-5.0%2FNOPB/363825?utm_adgroup=General&utm_source=google&utm_medium=cpc&utm_campaign=PMax%20Shopping_Product_Zombie%20SKUs&utm_term=&utm_content=General&utm_id=go_cmp-17815035045_adg-_ad-__dev-c_ext-_prd-363825_sig-CjwKCAiA5L2tBhBTEiwAdSxJX2jmx8jM-JlHhV04F58rlCzi0KZgwJl8jmcjRGNCM7uSaMTsq63izRoCBfYQAvD_BwE&gad_source=1&gclid=CjwKCAiA5L2tBhBTEiwAdSxJX2jmx8jM-JlHhV04F58rlCzi0KZgwJl8jmcjRGNCM7uSaMTsq63izRoCBfYQAvD_BwE"

Add-BomLink "F2"  "https://www.digikey.com/en/products/detail/stmicroelectronics/STM32F207VGT6TR/4357621" $null
Add-BomLink "F3"  $f3display $f3display
Add-BomLink "F4"  "https://www.digikey.com/en/products/detail/texas-instruments/LM2678S-3-3-NOPB/366918" $null
Add-BomLink "F5"  "https://www.digikey.com/en/products/detail/nichicon/UUD1H150MCL1GS/590040" $null
Add-BomLink "F6"  "https://www.digikey.com/en/products/detail/nichicon/UWP1HR47MCL1GB/2550802" $null
Add-BomLink "F7"  "https://www.digikey.com/en/products/detail/vishay-general-semiconductor-diodes-division/VS-6TQ045S-M3/5426222" $null
$ws.Range("F8").Value = "https://www.digikey.com/en/products/detail/bourns-inc/2205-H-RC/775358"
Add-BomLink "F8"  "https://www.digikey.com/en/products/detail/bourns-inc/2205-H-RC/775358" $null
Add-BomLink "F9"  "https://www.digikey.com/en/products/detail/nichicon/UCZ1J181MNJ1MS/5144110" $null
Add-BomLink "F10" "https://www.digikey.com/en/products/detail/murata-electronics/GRM2195C1H103JA01D/586788" $null

# ---------------------------------------------------------------------------
# 3. Append the new BOM rows (11-20) added for the updated schematic.
# ---------------------------------------------------------------------------
$ws.Range("A11").Value = "SN75ALS174ADWR"
$ws.Range("B11").Value = "Line Driver IC"
$ws.Range("C11").Value = 2
$ws.Range("F11").Value = "https://www.digikey.com/en/products/detail/texas-instruments/SN75ALS174ADWR/1593485"
Add-BomLink "F11" "https://www.digikey.com/en/products/detail/texas-instruments/SN75ALS174ADWR/1593485" $null

$ws.Range("A12").Value = "SN65LBC175AD"
$ws.Range("B12").Value = "Line Receiver IC"
$ws.Range("C12").Value = 4
Add-BomLink "F12" "https://www.digikey.com/en/products/detail/texas-instruments/SN65LBC175AD/380303" $null

$ws.Range("A13").Value = "PJ-202AH"
$ws.Range("B13").Value = "Barrel Jack Connector"
$ws.Range("C13").Value = 1
Add-BomLink "F13" "https://www.digikey.com/en/products/detail/cui-devices/PJ-202AH/408450" $null

$ws.Range("A14").Value = "OPA4205ADR"
$ws.Range("B14").Value = "General Op-Amp (4 channel)"
$ws.Range("C14").Value = 2
Add-BomLink "F14" "https://www.digikey.com/en/products/detail/texas-instruments/OPA4205APWR/17394950" $null

$ws.Range("A15").Value = "TPS63700DRCR"
$ws.Range("B15").Value = "Inverting Boost Converter"
$ws.Range("C15").Value = 1
Add-BomLink "F15" "https://www.digikey.com/en/products/detail/texas-instruments/TPS63700DRCR/1672393" $null

$ws.Range("A16").Value = "SL03-GS18"
$ws.Range("B16").Value = "Inv Boost Con. Diode"
$ws.Range("C16").Value = 1
Add-BomLink "F16" "https://www.digikey.com/en/products/detail/vishay-general-semiconductor-diodes-division/SL03-GS18/4871689" $null

$ws.Range("A17").Value = 7443551131
$ws.Range("B17").Value = "Inv Boost Con. Inductor"
$ws.Range("C17").Value = 1
Add-BomLink "F17" "https://www.digikey.com/en/products/detail/w%C3%BCrth-elektronik/7443551131/1638545" $null

$ws.Range("A18").Value = "3352T-1-203LF"
$ws.Range("B18").Value = "20k Potentiometer"
$ws.Range("C18").Value = 2
Add-BomLink "F18" "https://www.digikey.com/en/products/detail/bourns-inc/3352T-1-203LF/1088346" $null

$f19url = "https://www.digikey.com/en/products/detail/texas-instruments/SN75468DR/2255090?utm_adgroup=General&utm_source=google&utm_medium=cpc&utm_campaign=PMax%20Shopping_Product_Zombie%20SKUs&utm_term=&utm_content=General&utm_id=go_cmp-17815035045_adg-_ad-__dev-c_ext-_prd-2255090_sig-CjwKCAiAlJKuBhAdEiwAnZb7lY7edhjVnlVUhEyNawogcHzVo6bbfQ1LOtrzO4xh_eCL0cFOX98QUxoCbYMQAvD_BwE&gad_source=1&gclid=CjwKCAiAlJKuBhAdEiwAnZb7lY7edhjVnlVUhEyNawogcHzVo6bbfQ1LOtrzO4xh_eCL0cFOX98QUxoCbYMQAvD_BwE"
$ws.Range("A19").Value = "SN75468DR"
$ws.Range("B19").Value = "NPN Transistor Array"
$ws.Range("C19").Value = 2
Add-BomLink "F19" $f19url $f19url

$ws.Range("A20").Value = "ACSL-6400-00TE"
$ws.Range("B20").Value = "Optoisolator"
$ws.Range("C20").Value = 2
Add-BomLink "F20" "https://www.digikey.com/en/products/detail/broadcom-limited/ACSL-6400-00TE/825239" $null

# ---------------------------------------------------------------------------
# 4. Cosmetic sheet changes: wider Part-Description column and the cursor
#    left sitting on D12 like the source workbook.
# ---------------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 26
$ws.Range("D12").Select() | Out-Null
